# Add two new Hall of Fame ballots to the "ballots" sheet:
#   - Juan Vené      (row 37)
#   - Rob Biertempfel (row 38)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ballots")

# --- Row 37: Juan Vené ---
# Voted for: Todd Helton (F), Fred McGriff (J), Mike Mussina (K),
#            Andy Pettitte (M), Mariano Rivera (O)
$ws.Range("A37").Value = "Juan Vené"
$ws.Range("F37").Value = "x"
$ws.Range("J37").Value = "x"
$ws.Range("K37").Value = "x"
$ws.Range("M37").Value = "x"
$ws.Range("O37").Value = "x"
$ws.Range("AK37").Value = 5
$ws.Range("AL37").Value = "JuanVene.com"

# Copy the date formatting (numFmt) from the row above so we reuse the
# existing date style instead of creating a new one, then set the value.
$ws.Range("AM36").Copy()
$ws.Range("AM37").PasteSpecial(-4122)
$ws.Range("AM37").Value = 43447

# --- Row 38: Rob Biertempfel ---
# Voted for: Roy Halladay (E), Mike Mussina (K), Mariano Rivera (O),
#            Curt Schilling (Q), Larry Walker (V)
$ws.Range("A38").Value = "Rob Biertempfel"
$ws.Range("E38").Value = "x"
$ws.Range("K38").Value = "x"
$ws.Range("O38").Value = "x"
$ws.Range("Q38").Value = "x"
$ws.Range("V38").Value = "x"
$ws.Range("AK38").Value = 5
$ws.Range("AL38").Value = "twitter"

$ws.Range("AM36").Copy()
$ws.Range("AM38").PasteSpecial(-4122)
$ws.Range("AM38").Value = 43447

# Match the author's final selection in the sheet view.
$ws.Range("B35").Select()
